$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ E=3; F=1; G=36.899643; H=110.698929; I=0.7238945645409351; J=0.7238945645409351; K=3; L=1; M=6.346253666666667; N=19.038761; O=0.9446330608455225; P=0.9446330608455226; Q=234.174494687441; R=2107.570452186969; S=0.6838147382317402; T=0.6838147382317402 }
    3 = @{ E=3; F=1; G=36.899643; H=110.698929; I=0.7238945645409351; J=0.7238945645409351; K=3; L=1; M=0.3719673333333333; N=1.115902; O=0.05536693915447755; P=0.05536693915447755; Q=13.725461807662; R=123.529156268958; S=0.04007982630919498; T=0.04007982630919498 }
    4 = @{ E=3; F=1; G=3.374819; H=10.124457; I=0.0662069584361419; J=0.0662069584361419; K=3; L=1; M=6.346253666666667; N=19.038761; O=0.9446330608455225; P=0.9446330608455226; Q=21.41745745308633; R=192.757117077777; S=0.062541281796805; T=0.06254128179680501 }
    5 = @{ E=3; F=1; G=3.374819; H=10.124457; I=0.0662069584361419; J=0.0662069584361419; K=3; L=1; M=0.3719673333333333; N=1.115902; O=0.05536693915447755; P=0.05536693915447755; Q=1.255322423912667; R=11.297901815214; S=0.003665676639336893; T=0.003665676639336893 }
    6 = @{ E=3; F=1; G=10.699319; H=32.097957; I=0.2098984770229228; J=0.2098984770229228; K=3; L=1; M=6.346253666666667; N=19.038761; O=0.9446330608455225; P=0.9446330608455226; Q=67.90059243458634; R=611.1053319112771; S=0.1982770408169772; T=0.1982770408169772 }
    7 = @{ E=3; F=1; G=10.699319; H=32.097957; I=0.2098984770229228; J=0.2098984770229228; K=3; L=1; M=0.3719673333333333; N=1.115902; O=0.05536693915447755; P=0.05536693915447755; Q=3.979797156912667; R=35.818174412214; S=0.01162143620594567; T=0.01162143620594567 }
}

foreach ($rowNum in $data.Keys) {
    $rowVals = $data[$rowNum]
    foreach ($col in $rowVals.Keys) {
        $ws.Range("$col$rowNum").Value = $rowVals[$col]
    }
}
